$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "602"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "vending"
$ws.Range("C3").Value = 30
